$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.472.85'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.877.99'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.023'
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.71'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.020'
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5122'
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3955'
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08480'
$ws.Range("E9").Value = '  +1.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.258'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.873.01'
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.47'
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.251'
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.024'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001112'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.16'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06799'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.71'
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.019'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.944'
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.517.54'
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.17'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.291'
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.093.13'
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.93'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.81'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.350'
$ws.Range("E28").Value = '  -4.44%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '127.25'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1057'
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.044'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.775'
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.645'
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.02432'
$ws.Range("E34").Value = '  -1.39%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06485'
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2177'
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.866'
$ws.Range("E37").Value = '  -6.95%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.268'
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.187'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6390'
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.002'
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.22'
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("B43").Value = 'Decentraland'
$ws.Range("C43").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6051'
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.97'
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.724'
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.994'
$ws.Range("E46").Value = '  -1.13%  '
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.208'
$ws.Range("E47").Value = '  -6.27%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.09'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.207'
$ws.Range("E49").Value = '  -2.82%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06850'
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '76.69'
$ws.Range("E51").Value = '  -1.90%  '
